# Updated cryptos list on Sat Apr 13 04:38:05 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# swaps the BitcoinCash / WrappedEther rows (17 and 18) to reflect the new order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "1.00", "152.90") are force-
# formatted as Text first so Excel keeps them as the literal string instead of
# silently normalising them into a numeric value (which would drop trailing zeros).

$ws.Range("D2").Value = '67.282.83'
$ws.Range("E2").Value = '  -5.07%  '
$ws.Range("D3").Value = '3.233.87'
$ws.Range("E3").Value = '  -8.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.91'
$ws.Range("E5").Value = '  -5.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.90'
$ws.Range("E6").Value = '  -12.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '3.227.35'
$ws.Range("E8").Value = '  -8.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("E9").Value = '  -10.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.173'
$ws.Range("E10").Value = '  -12.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.75'
$ws.Range("E11").Value = '  -6.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.505'
$ws.Range("E12").Value = '  -14.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.41'
$ws.Range("E13").Value = '  -17.80%  '
$ws.Range("E14").Value = '  -11.69%  '
$ws.Range("D15").Value = '3.754.02'
$ws.Range("E15").Value = '  -8.72%  '
$ws.Range("D16").Value = '67.221.44'
$ws.Range("E16").Value = '  -5.17%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.233.58'
$ws.Range("E17").Value = '  -8.04%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '541.74'
$ws.Range("E18").Value = '  -11.58%  '
$ws.Range("E19").Value = '  -5.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.18'
$ws.Range("E20").Value = '  -15.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.20'
$ws.Range("E21").Value = '  -14.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.762'
$ws.Range("E22").Value = '  -14.39%  '
$ws.Range("E23").Value = '  -14.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.69'
$ws.Range("E24").Value = '  -12.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.53'
$ws.Range("E25").Value = '  -14.09%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -16.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.11'
$ws.Range("E28").Value = '  -11.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '29.47'
$ws.Range("E29").Value = '  -13.05%  '
$ws.Range("E30").Value = '  -17.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.69'
$ws.Range("E31").Value = '  -11.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.15'
$ws.Range("E32").Value = '  -11.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '543.37'
$ws.Range("E33").Value = '  -9.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.57'
$ws.Range("E34").Value = '  -19.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.75'
$ws.Range("E35").Value = '  -16.43%  '
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0448'
$ws.Range("E37").Value = '  -5.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.70'
$ws.Range("E38").Value = '  -5.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0849'
$ws.Range("E39").Value = '  -15.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.20'
$ws.Range("E40").Value = '  -15.30%  '
$ws.Range("E41").Value = '  -13.05%  '
$ws.Range("D42").Value = '2.927.54'
$ws.Range("E42").Value = '  -13.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.61'
$ws.Range("E43").Value = '  -26.05%  '
$ws.Range("D44").Value = '0.0₃0587'
$ws.Range("E44").Value = '  -20.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.261'
$ws.Range("E45").Value = '  -16.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.39'
$ws.Range("E46").Value = '  -20.32%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.17'
$ws.Range("E48").Value = '  -18.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.12'
$ws.Range("E49").Value = '  -17.52%  '
$ws.Range("E50").Value = '  -13.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '124.37'
$ws.Range("E51").Value = '  -7.03%  '
